$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.301.32"
$ws.Range("E2").Value = "  +2.59%  "

$ws.Range("D3").Value = "2.587.24"
$ws.Range("E3").Value = "  +10.87%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.42%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.09"
$ws.Range("E5").Value = "  +2.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.31"
$ws.Range("E6").Value = "  +6.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.610"
$ws.Range("E7").Value = "  +7.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.578"
$ws.Range("E9").Value = "  +13.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.20"
$ws.Range("E10").Value = "  +13.33%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0837"
$ws.Range("E11").Value = "  +5.97%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.02"
$ws.Range("E12").Value = "  +12.91%  "

$ws.Range("D13").Value = "2.971.40"
$ws.Range("E13").Value = "  +10.45%  "

$ws.Range("E14").Value = "  +2.93%  "

$ws.Range("D15").Value = "2.600.05"
$ws.Range("E15").Value = "  +11.38%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.906"
$ws.Range("E16").Value = "  +12.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "15.22"
$ws.Range("E17").Value = "  +11.59%  "

$ws.Range("D18").Value = "47.225.69"
$ws.Range("E18").Value = "  +2.47%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.15"
$ws.Range("E19").Value = "  +11.80%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000100"
$ws.Range("E20").Value = "  +3.95%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.70"
$ws.Range("E21").Value = "  +12.54%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.53"
$ws.Range("E22").Value = "  +6.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "258.38"
$ws.Range("E23").Value = "  +5.91%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.00"
$ws.Range("E24").Value = "  +6.78%  "

$ws.Range("E25").Value = "  +12.40%  "

$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.70"
$ws.Range("E27").Value = "  +23.78%  "

$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "42.33"
$ws.Range("E28").Value = "  +5.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.56"
$ws.Range("E29").Value = "  +9.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.28"
$ws.Range("E30").Value = "  +2.83%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.88"
$ws.Range("E31").Value = "  +7.99%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.04"
$ws.Range("E32").Value = "  +11.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.97"
$ws.Range("E33").Value = "  +5.65%  "

$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.25"
$ws.Range("E34").Value = "  +25.60%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0851"
$ws.Range("E35").Value = "  +10.36%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "150.54"
$ws.Range("E36").Value = "  +4.13%  "

$ws.Range("E37").Value = "  +9.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.122"
$ws.Range("E38").Value = "  +4.89%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.64"
$ws.Range("E39").Value = "  +9.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.31"
$ws.Range("E40").Value = "  +11.48%  "

$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.63"
$ws.Range("E41").Value = "  +14.18%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0330"
$ws.Range("E42").Value = "  +10.84%  "

$ws.Range("D43").Value = "2.011.48"
$ws.Range("E43").Value = "  +8.26%  "

$ws.Range("E44").Value = "  -0.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "95.23"
$ws.Range("E45").Value = "  +5.61%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.95"
$ws.Range("E46").Value = "  +39.53%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.93"
$ws.Range("E47").Value = "  +6.90%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.203"
$ws.Range("E48").Value = "  +10.15%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "108.76"
$ws.Range("E49").Value = "  +13.72%  "

$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.97"
$ws.Range("E50").Value = "  +13.28%  "

$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "76.05"
$ws.Range("E51").Value = "  +8.70%  "
